$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B23").Value = "Assign Alfresco Folder"
$ws.Range("C23").Value = "container?.folder?.cmisFolderId == null"
$ws.Range("D23").Value = "setEcmFolderPath, '/Sites/acm/documentLibrary/Timesheets/' + dateFormat('yyyyMMdd') + '_' + `$acmTimesheet.getId()"

# Column D width change
$ws.Columns.Item(4).ColumnWidth = 112.140625

# Scroll / selection changes
$ws.Range("D37").Select()
$excel.ActiveWindow.ScrollRow = 16
